$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.531.99"
$ws.Cells.Item(2, 5).Value = "  +1.52%  "

$ws.Cells.Item(3, 4).Value = "1.825.59"
$ws.Cells.Item(3, 5).Value = "  +1.40%  "

$ws.Cells.Item(4, 4).Value = "1.003"
$ws.Cells.Item(4, 5).Value = "  +0.14%  "

$ws.Cells.Item(5, 4).Value = "316.13"
$ws.Cells.Item(5, 5).Value = "  -0.08%  "

$ws.Cells.Item(6, 4).Value = "1.003"
$ws.Cells.Item(6, 5).Value = "  +0.14%  "

$ws.Cells.Item(7, 4).Value = "0.5304"
$ws.Cells.Item(7, 5).Value = "  -2.73%  "

$ws.Cells.Item(8, 4).Value = "0.3937"
$ws.Cells.Item(8, 5).Value = "  +4.12%  "

$ws.Cells.Item(9, 4).Value = "0.07726"
$ws.Cells.Item(9, 5).Value = "  +3.56%  "

$ws.Cells.Item(10, 4).Value = "42.01"
$ws.Cells.Item(10, 5).Value = "  +0.13%  "

$ws.Cells.Item(11, 4).Value = "1.113"
$ws.Cells.Item(11, 5).Value = "  +1.86%  "

$ws.Cells.Item(12, 4).Value = "21.05"
$ws.Cells.Item(12, 5).Value = "  +3.08%  "

$ws.Cells.Item(13, 4).Value = "6.295"
$ws.Cells.Item(13, 5).Value = "  +1.55%  "

$ws.Cells.Item(14, 2).Value = "BinanceUSD"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(14, 4).Value = "1.003"
$ws.Cells.Item(14, 5).Value = "  +0.20%  "

$ws.Cells.Item(15, 2).Value = "Chainlink"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(15, 4).Value = "7.559"
$ws.Cells.Item(15, 5).Value = "  +2.95%  "

$ws.Cells.Item(16, 4).Value = "1.822.88"
$ws.Cells.Item(16, 5).Value = "  +1.38%  "

$ws.Cells.Item(17, 4).Value = "92.97"
$ws.Cells.Item(17, 5).Value = "  +3.82%  "

$ws.Cells.Item(18, 4).Value = "'0.00001082"
$ws.Cells.Item(18, 5).Value = "  +1.72%  "

$ws.Cells.Item(19, 4).Value = "0.06607"
$ws.Cells.Item(19, 5).Value = "  +0.98%  "

$ws.Cells.Item(20, 4).Value = "'17.70"
$ws.Cells.Item(20, 5).Value = "  +1.60%  "

$ws.Cells.Item(21, 4).Value = "1.002"
$ws.Cells.Item(21, 5).Value = "  +0.15%  "

$ws.Cells.Item(22, 4).Value = "6.074"
$ws.Cells.Item(22, 5).Value = "  +2.34%  "

$ws.Cells.Item(23, 4).Value = "28.537.82"
$ws.Cells.Item(23, 5).Value = "  +1.41%  "

$ws.Cells.Item(24, 4).Value = "11.13"
$ws.Cells.Item(24, 5).Value = "  -0.52%  "

$ws.Cells.Item(25, 4).Value = "2.239"
$ws.Cells.Item(25, 5).Value = "  +7.34%  "

$ws.Cells.Item(26, 4).Value = "20.68"
$ws.Cells.Item(26, 5).Value = "  +1.28%  "

$ws.Cells.Item(27, 4).Value = "156.64"
$ws.Cells.Item(27, 5).Value = "  +0.70%  "

$ws.Cells.Item(28, 4).Value = "2.035.73"
$ws.Cells.Item(28, 5).Value = "  +1.50%  "

$ws.Cells.Item(29, 4).Value = "2.413"
$ws.Cells.Item(29, 5).Value = "  +3.96%  "

$ws.Cells.Item(30, 4).Value = "125.18"
$ws.Cells.Item(30, 5).Value = "  +2.83%  "

$ws.Cells.Item(31, 4).Value = "1.137"
$ws.Cells.Item(31, 5).Value = "  +1.91%  "

$ws.Cells.Item(32, 4).Value = "0.1119"
$ws.Cells.Item(32, 5).Value = "  +0.44%  "

$ws.Cells.Item(33, 4).Value = "5.713"
$ws.Cells.Item(33, 5).Value = "  +2.74%  "

$ws.Cells.Item(34, 4).Value = "3.658"
$ws.Cells.Item(34, 5).Value = "  -0.29%  "

$ws.Cells.Item(35, 4).Value = "0.07281"
$ws.Cells.Item(35, 5).Value = "  +5.34%  "

$ws.Cells.Item(36, 4).Value = "0.2251"
$ws.Cells.Item(36, 5).Value = "  +1.33%  "

$ws.Cells.Item(37, 4).Value = "0.02347"
$ws.Cells.Item(37, 5).Value = "  +2.49%  "

$ws.Cells.Item(38, 4).Value = "8.882"
$ws.Cells.Item(38, 5).Value = "  +5.03%  "

$ws.Cells.Item(39, 4).Value = "5.163"
$ws.Cells.Item(39, 5).Value = "  +1.47%  "

$ws.Cells.Item(40, 4).Value = "11.34"
$ws.Cells.Item(40, 5).Value = "  +1.51%  "

$ws.Cells.Item(41, 4).Value = "0.6264"
$ws.Cells.Item(41, 5).Value = "  +1.76%  "

$ws.Cells.Item(42, 4).Value = "1.191"
$ws.Cells.Item(42, 5).Value = "  +1.54%  "

$ws.Cells.Item(43, 4).Value = "1.002"
$ws.Cells.Item(43, 5).Value = "  +0.14%  "

$ws.Cells.Item(44, 4).Value = "1.396"
$ws.Cells.Item(44, 5).Value = "  -1.79%  "

$ws.Cells.Item(45, 4).Value = "13.45"
$ws.Cells.Item(45, 5).Value = "  +1.28%  "

$ws.Cells.Item(46, 4).Value = "0.5916"
$ws.Cells.Item(46, 5).Value = "  +3.05%  "

$ws.Cells.Item(47, 5).Value = "  +1.10%  "

$ws.Cells.Item(48, 4).Value = "125.47"
$ws.Cells.Item(48, 5).Value = "  +0.61%  "

$ws.Cells.Item(49, 4).Value = "1.991"
$ws.Cells.Item(49, 5).Value = "  +3.85%  "

$ws.Cells.Item(50, 4).Value = "'1.190"
$ws.Cells.Item(50, 5).Value = "  +0.63%  "

$ws.Cells.Item(51, 4).Value = "0.06951"
$ws.Cells.Item(51, 5).Value = "  +2.02%  "
